# Add a new "2022" column (K) to the right of the existing 2015-2021 data
# (column J), mirroring column J's formatting, then update the sheet's
# selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column J's formatting (rows 4-14, the data block) into the new
# column K so the new cells pick up the same styles (number format,
# borders, fonts, etc.) that the rest of the table uses.
$ws.Range("J4:J14").Copy()
$ws.Range("K4:K14").PasteSpecial(-4122)   # xlPasteFormats
[void]($excel.CutCopyMode = $false)

# Write the 2022 values into the new column.
$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.3
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Match the author's final active-cell selection.
[void]$ws.Range("L7").Select()
